# Auto-generated edit script for before.xlsx -> target state
# Source: diff of data/recommandations.xlsx (BRVM auto-update GitHub Action)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- "Recommandations" sheet --------------------------------------------
# The underlying BRVM scrape re-ranked / re-derived this table, so rows
# 2-44 are rewritten in full (titles, day counts, variations, and the
# recommendation/strategy labels), and the now-unused trailing rows
# 45-48 are removed so the sheet shrinks from A1:G48 to A1:G44.

# Row 2: BRVM - SERVICES PUBLICS
$ws1.Range("A2").Value = 'BRVM - SERVICES PUBLICS'
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = 3263.03
$ws1.Range("E2").Value = 99.7
$ws1.Range("F2").Value = '🟡 Observer'
$ws1.Range("G2").Value = '➖ Neutre'

# Row 3: UNIWAX CI
$ws1.Range("A3").Value = 'UNIWAX CI'
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 2850
$ws1.Range("E3").Value = 730
$ws1.Range("F3").Value = '🟡 Observer'
$ws1.Range("G3").Value = '➖ Neutre'

# Row 4: CFAO MOTORS CI
$ws1.Range("A4").Value = 'CFAO MOTORS CI'
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 2670
$ws1.Range("E4").Value = 680
$ws1.Range("F4").Value = '🟡 Observer'
$ws1.Range("G4").Value = '➖ Neutre'

# Row 5: BRVM - AUTRES SECTEURS
$ws1.Range("A5").Value = 'BRVM - AUTRES SECTEURS'
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 2528.56
$ws1.Range("E5").Value = 633.23
$ws1.Range("F5").Value = '🟡 Observer'
$ws1.Range("G5").Value = '➖ Neutre'

# Row 6: NEI-CEDA CI
$ws1.Range("A6").Value = 'NEI-CEDA CI'
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 2395
$ws1.Range("E6").Value = 600
$ws1.Range("F6").Value = '🟡 Observer'
$ws1.Range("G6").Value = '➖ Neutre'

# Row 7: SETAO CI
$ws1.Range("A7").Value = 'SETAO CI'
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 2365
$ws1.Range("E7").Value = 590
$ws1.Range("F7").Value = '🟡 Observer'
$ws1.Range("G7").Value = '➖ Neutre'

# Row 8: AIR LIQUIDE CI
$ws1.Range("A8").Value = 'AIR LIQUIDE CI'
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 2275
$ws1.Range("E8").Value = 590
$ws1.Range("F8").Value = '🟡 Observer'
$ws1.Range("G8").Value = '➖ Neutre'

# Row 9: BRVM - DISTRIBUTION
$ws1.Range("A9").Value = 'BRVM - DISTRIBUTION'
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 4
$ws1.Range("D9").Value = 1489.92
$ws1.Range("E9").Value = 372.07
$ws1.Range("F9").Value = '🟡 Observer'
$ws1.Range("G9").Value = '➖ Neutre'

# Row 10: BRVM - TRANSPORT
$ws1.Range("A10").Value = 'BRVM - TRANSPORT'
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 4
$ws1.Range("D10").Value = 1488.23
$ws1.Range("E10").Value = 375.73
$ws1.Range("F10").Value = '🟡 Observer'
$ws1.Range("G10").Value = '➖ Neutre'

# Row 11: BRVM - AGRICULTURE
$ws1.Range("A11").Value = 'BRVM - AGRICULTURE'
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 4
$ws1.Range("D11").Value = 1312.03
$ws1.Range("E11").Value = 331.16
$ws1.Range("F11").Value = '🟡 Observer'
$ws1.Range("G11").Value = '➖ Neutre'

# Row 12: BRVM - INDUSTRIE
$ws1.Range("A12").Value = 'BRVM - INDUSTRIE'
$ws1.Range("B12").Value = 0
$ws1.Range("C12").Value = 4
$ws1.Range("D12").Value = 1068.97
$ws1.Range("E12").Value = 268.78
$ws1.Range("F12").Value = '🟡 Observer'
$ws1.Range("G12").Value = '➖ Neutre'

# Row 13: BRVM - CONSOMMATION DE BASE
$ws1.Range("A13").Value = 'BRVM - CONSOMMATION DE BASE'
$ws1.Range("B13").Value = 0
$ws1.Range("C13").Value = 4
$ws1.Range("D13").Value = 882.08
$ws1.Range("E13").Value = 221.93
$ws1.Range("F13").Value = '🟡 Observer'
$ws1.Range("G13").Value = '➖ Neutre'

# Row 14: BRVM-PRINCIPAL
$ws1.Range("A14").Value = 'BRVM-PRINCIPAL'
$ws1.Range("B14").Value = 0
$ws1.Range("C14").Value = 4
$ws1.Range("D14").Value = 769.4400000000001
$ws1.Range("E14").Value = 193.09
$ws1.Range("F14").Value = '🟡 Observer'
$ws1.Range("G14").Value = '➖ Neutre'

# Row 15: BRVM - INDUSTRIELS
$ws1.Range("A15").Value = 'BRVM - INDUSTRIELS'
$ws1.Range("B15").Value = 0
$ws1.Range("C15").Value = 4
$ws1.Range("D15").Value = 565
$ws1.Range("E15").Value = 142.49
$ws1.Range("F15").Value = '🟡 Observer'
$ws1.Range("G15").Value = '➖ Neutre'

# Row 16: BRVM-PRESTIGE
$ws1.Range("A16").Value = 'BRVM-PRESTIGE'
$ws1.Range("B16").Value = 0
$ws1.Range("C16").Value = 4
$ws1.Range("D16").Value = 527.29
$ws1.Range("E16").Value = 132.12
$ws1.Range("F16").Value = '🟡 Observer'
$ws1.Range("G16").Value = '➖ Neutre'

# Row 17: BRVM - FINANCES
$ws1.Range("A17").Value = 'BRVM - FINANCES'
$ws1.Range("B17").Value = 0
$ws1.Range("C17").Value = 4
$ws1.Range("D17").Value = 497.65
$ws1.Range("E17").Value = 124.83
$ws1.Range("F17").Value = '🟡 Observer'
$ws1.Range("G17").Value = '➖ Neutre'

# Row 18: BRVM - SERVICES FINANCIERS
$ws1.Range("A18").Value = 'BRVM - SERVICES FINANCIERS'
$ws1.Range("B18").Value = 0
$ws1.Range("C18").Value = 4
$ws1.Range("D18").Value = 489.09
$ws1.Range("E18").Value = 122.68
$ws1.Range("F18").Value = '🟡 Observer'
$ws1.Range("G18").Value = '➖ Neutre'

# Row 19: BRVM - ENERGIE
$ws1.Range("A19").Value = 'BRVM - ENERGIE'
$ws1.Range("B19").Value = 0
$ws1.Range("C19").Value = 4
$ws1.Range("D19").Value = 442.28
$ws1.Range("E19").Value = 111.26
$ws1.Range("F19").Value = '🟡 Observer'
$ws1.Range("G19").Value = '➖ Neutre'

# Row 20: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Range("A20").Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws1.Range("B20").Value = 0
$ws1.Range("C20").Value = 4
$ws1.Range("D20").Value = 430.57
$ws1.Range("E20").Value = 106.95
$ws1.Range("F20").Value = '🟡 Observer'
$ws1.Range("G20").Value = '➖ Neutre'

# Row 21: BRVM - TELECOMMUNICATIONS
$ws1.Range("A21").Value = 'BRVM - TELECOMMUNICATIONS'
$ws1.Range("B21").Value = 0
$ws1.Range("C21").Value = 4
$ws1.Range("D21").Value = 374.09
$ws1.Range("E21").Value = 93.61
$ws1.Range("F21").Value = '🟡 Observer'
$ws1.Range("G21").Value = '➖ Neutre'

# Row 22: SAFCA CI (SAFC)
$ws1.Range("A22").Value = 'SAFCA CI (SAFC)'
$ws1.Range("B22").Value = 4
$ws1.Range("C22").Value = 0
$ws1.Range("D22").Value = 24.81
$ws1.Range("E22").Value = 3.27
$ws1.Range("F22").Value = '🟢 Achat'
$ws1.Range("G22").Value = '✅ Renforcer'

# Row 23: SUCRIVOIRE (SCRC)
$ws1.Range("A23").Value = 'SUCRIVOIRE (SCRC)'
$ws1.Range("B23").Value = 3
$ws1.Range("C23").Value = 0
$ws1.Range("D23").Value = 16.95
$ws1.Range("E23").Value = 2.95
$ws1.Range("F23").Value = '🟢 Achat'
$ws1.Range("G23").Value = '✅ Renforcer'

# Row 24: UNIWAX CI (UNXC)
$ws1.Range("A24").Value = 'UNIWAX CI (UNXC)'
$ws1.Range("B24").Value = 3
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = 14.21
$ws1.Range("E24").Value = -4.05
$ws1.Range("F24").Value = '🟢 Achat'
$ws1.Range("G24").Value = '✅ Renforcer'

# Row 25: BERNABE CI (BNBC)
$ws1.Range("A25").Value = 'BERNABE CI (BNBC)'
$ws1.Range("B25").Value = 3
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = 13.74
$ws1.Range("E25").Value = -2.81
$ws1.Range("F25").Value = '🟢 Achat'
$ws1.Range("G25").Value = '✅ Renforcer'

# Row 26: ORAGROUP TOGO (ORGT)
$ws1.Range("A26").Value = 'ORAGROUP TOGO (ORGT)'
$ws1.Range("B26").Value = 2
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = 7.45
$ws1.Range("E26").Value = 4.39
$ws1.Range("F26").Value = '🟡 Observer'
$ws1.Range("G26").Value = '👀 À surveiller'

# Row 27: AIR LIQUIDE CI (SIVC)
$ws1.Range("A27").Value = 'AIR LIQUIDE CI (SIVC)'
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 0
$ws1.Range("D27").Value = 7.27
$ws1.Range("E27").Value = 7.27
$ws1.Range("F27").Value = '🟡 Observer'
$ws1.Range("G27").Value = '➖ Neutre'

# Row 28: SETAO CI (STAC)
$ws1.Range("A28").Value = 'SETAO CI (STAC)'
$ws1.Range("B28").Value = 1
$ws1.Range("C28").Value = 1
$ws1.Range("D28").Value = 4.24
$ws1.Range("E28").Value = 5.08
$ws1.Range("F28").Value = '🟡 Observer'
$ws1.Range("G28").Value = '👀 À surveiller'

# Row 29: SERVAIR ABIDJAN CI (ABJC)
$ws1.Range("A29").Value = 'SERVAIR ABIDJAN CI (ABJC)'
$ws1.Range("B29").Value = 1
$ws1.Range("C29").Value = 0
$ws1.Range("D29").Value = 2.56
$ws1.Range("E29").Value = 2.56
$ws1.Range("F29").Value = '🟡 Observer'
$ws1.Range("G29").Value = '➖ Neutre'

# Row 30: TOTAL
$ws1.Range("A30").Value = 'TOTAL'
$ws1.Range("B30").Value = 0
$ws1.Range("C30").Value = 4
$ws1.Range("D30").Value = 0
$ws1.Range("E30").Value = 0
$ws1.Range("F30").Value = '🟡 Observer'
$ws1.Range("G30").Value = '➖ Neutre'

# Row 31: TOTALENERGIES MARKETING SN (TTLS)
$ws1.Range("A31").Value = 'TOTALENERGIES MARKETING SN (TTLS)'
$ws1.Range("B31").Value = 0
$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = -0.2
$ws1.Range("E31").Value = -0.2
$ws1.Range("F31").Value = '🟡 Observer'
$ws1.Range("G31").Value = '➖ Neutre'

# Row 32: VIVO ENERGY CI (SHEC)
$ws1.Range("A32").Value = 'VIVO ENERGY CI (SHEC)'
$ws1.Range("B32").Value = 1
$ws1.Range("C32").Value = 1
$ws1.Range("D32").Value = -0.36
$ws1.Range("E32").Value = 2.17
$ws1.Range("F32").Value = '🟡 Observer'
$ws1.Range("G32").Value = '👀 À surveiller'

# Row 33: ECOBANK COTE D''IVOIRE (ECOC)
$ws1.Range("A33").Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Range("B33").Value = 0
$ws1.Range("C33").Value = 1
$ws1.Range("D33").Value = -0.42
$ws1.Range("E33").Value = -0.42
$ws1.Range("F33").Value = '🟡 Observer'
$ws1.Range("G33").Value = '➖ Neutre'

# Row 34: LOTERIE NATIONALE DU BENIN (LNBB)
$ws1.Range("A34").Value = 'LOTERIE NATIONALE DU BENIN (LNBB)'
$ws1.Range("B34").Value = 0
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = -0.55
$ws1.Range("E34").Value = -0.55
$ws1.Range("F34").Value = '🟡 Observer'
$ws1.Range("G34").Value = '➖ Neutre'

# Row 35: SOCIETE IVOIRIENNE DE BANQUE  (SIBC)
$ws1.Range("A35").Value = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$ws1.Range("B35").Value = 0
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = -0.82
$ws1.Range("E35").Value = -0.82
$ws1.Range("F35").Value = '🟡 Observer'
$ws1.Range("G35").Value = '➖ Neutre'

# Row 36: NSIA BANQUE COTE D'IVOIRE (NSBC)
$ws1.Range("A36").Value = 'NSIA BANQUE COTE D''IVOIRE (NSBC)'
$ws1.Range("B36").Value = 0
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = -0.89
$ws1.Range("E36").Value = -0.89
$ws1.Range("F36").Value = '🟡 Observer'
$ws1.Range("G36").Value = '➖ Neutre'

# Row 37: SOLIBRA CI (SLBC)
$ws1.Range("A37").Value = 'SOLIBRA CI (SLBC)'
$ws1.Range("B37").Value = 1
$ws1.Range("C37").Value = 1
$ws1.Range("D37").Value = -1.23
$ws1.Range("E37").Value = 3.83
$ws1.Range("F37").Value = '🟡 Observer'
$ws1.Range("G37").Value = '👀 À surveiller'

# Row 38: ORANGE COTE D'IVOIRE (ORAC)
$ws1.Range("A38").Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$ws1.Range("B38").Value = 0
$ws1.Range("C38").Value = 1
$ws1.Range("D38").Value = -1.7
$ws1.Range("E38").Value = -1.7
$ws1.Range("F38").Value = '🟡 Observer'
$ws1.Range("G38").Value = '➖ Neutre'

# Row 39: CIE CI (CIEC)
$ws1.Range("A39").Value = 'CIE CI (CIEC)'
$ws1.Range("B39").Value = 0
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = -1.86
$ws1.Range("E39").Value = -1.86
$ws1.Range("F39").Value = '🟡 Observer'
$ws1.Range("G39").Value = '➖ Neutre'

# Row 40: CFAO MOTORS CI (CFAC)
$ws1.Range("A40").Value = 'CFAO MOTORS CI (CFAC)'
$ws1.Range("B40").Value = 0
$ws1.Range("C40").Value = 1
$ws1.Range("D40").Value = -2.21
$ws1.Range("E40").Value = -2.21
$ws1.Range("F40").Value = '🟡 Observer'
$ws1.Range("G40").Value = '➖ Neutre'

# Row 41: SICABLE CI (CABC)
$ws1.Range("A41").Value = 'SICABLE CI (CABC)'
$ws1.Range("B41").Value = 0
$ws1.Range("C41").Value = 1
$ws1.Range("D41").Value = -2.69
$ws1.Range("E41").Value = -2.69
$ws1.Range("F41").Value = '🟡 Observer'
$ws1.Range("G41").Value = '➖ Neutre'

# Row 42: BICI CI (BICC)
$ws1.Range("A42").Value = 'BICI CI (BICC)'
$ws1.Range("B42").Value = 0
$ws1.Range("C42").Value = 2
$ws1.Range("D42").Value = -2.73
$ws1.Range("E42").Value = -0.92
$ws1.Range("F42").Value = '🟡 Observer'
$ws1.Range("G42").Value = '➖ Neutre'

# Row 43: ONATEL BF (ONTBF)
$ws1.Range("A43").Value = 'ONATEL BF (ONTBF)'
$ws1.Range("B43").Value = 0
$ws1.Range("C43").Value = 2
$ws1.Range("D43").Value = -2.91
$ws1.Range("E43").Value = -1.67
$ws1.Range("F43").Value = '🟡 Observer'
$ws1.Range("G43").Value = '➖ Neutre'

# Row 44: NEI-CEDA CI (NEIC)
$ws1.Range("A44").Value = 'NEI-CEDA CI (NEIC)'
$ws1.Range("B44").Value = 0
$ws1.Range("C44").Value = 1
$ws1.Range("D44").Value = -6.67
$ws1.Range("E44").Value = -6.67
$ws1.Range("F44").Value = '🟡 Observer'
$ws1.Range("G44").Value = '➖ Neutre'

# Drop the 4 rows that no longer exist in the refreshed dataset so the
# sheet dimension goes back to A1:G44.
$ws1.Rows("45:48").Delete()

# --- "Top_YTD" sheet -----------------------------------------------------
# Only the YTD progress figures (column B) are refreshed; titles (col A)
# are unchanged.

$ws2.Range("B2").Value = 7056999.35  # BRVM - SERVICES PUBLICS
$ws2.Range("B3").Value = 434491.53  # UNIWAX CI
$ws2.Range("B4").Value = 346731.88  # CFAO MOTORS CI
$ws2.Range("B5").Value = 287226.11  # BRVM - AUTRES SECTEURS
$ws2.Range("B6").Value = 238285  # NEI-CEDA CI
$ws2.Range("B7").Value = 228165.8  # SETAO CI
$ws2.Range("B8").Value = 199504.92  # AIR LIQUIDE CI
$ws2.Range("B9").Value = 49734.44  # BRVM - DISTRIBUTION
$ws2.Range("B10").Value = 49552.55  # BRVM - TRANSPORT
$ws2.Range("B11").Value = 33456.04  # BRVM - AGRICULTURE
